$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the English keywords back to the original Chinese keywords.
$ws.Range("B2").Value = "争吵"

$ws.Range("B3").Value = "考验"
$ws.Range("C3").Value = "逃避"

$ws.Range("B4").Value = "界线感"
$ws.Range("C4").Value = "放下戒备"
$ws.Range("D4").Value = "信任"

$ws.Range("B5").Value = "夏天"
$ws.Range("C5").Value = "失望"
$ws.Range("D5").Value = "放弃"
$ws.Range("E5").Value = "分手"

# Update the active selection to match the saved view state.
$ws.Range("E8").Select()
